$d = $word.ActiveDocument

# --- Add new table row: n = 9, a(9) = 10599739324 ---
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "9"
$newRow.Cells.Item(2).Range.Text = "10599739324"

# --- Update the cached date field result in the default footer ---
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ftr.Range.Find.Execute("04.12.2016", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "05.12.2016", 2)
